# FLRenamingConfigFile.xlsx - "Adjusted the ConfigFile for testing"
#
# The "Config" sheet's two recipient-list rows (row 10 = RecipientTo,
# row 11 = RecipientCC) are both collapsed down to the single tester
# address "lester.rollan@lexisnexisrisk.com" so that test runs of the
# automation only email the one address instead of the full
# production distribution lists.

$wb = $excel.ActiveWorkbook

$configSheet = $wb.Worksheets.Item("Config")

$testerEmail = "lester.rollan@lexisnexisrisk.com"

# RecipientTo (row 10) and RecipientCC (row 11) -> single tester address
$configSheet.Range("B10").Value = $testerEmail
$configSheet.Range("B11").Value = $testerEmail

# Reflect where the author was last working on the sheet: they had
# scrolled/selected down near the bottom of the filled-in rows.
$configSheet.Activate()
$configSheet.Range("B15").Select()
